$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.842093
$ws.Range("H2").Value = 11.526279
$ws.Range("I2").Value = 0.8773182083458525
$ws.Range("J2").Value = 0.886145240040679
$ws.Range("M2").Value = 1.009860666666667
$ws.Range("N2").Value = 3.029582
$ws.Range("O2").Value = 0.01353413605720072
$ws.Range("P2").Value = 0.01542521070970148
$ws.Range("Q2").Value = 3.879978598375334
$ws.Range("R2").Value = 34.919807385378
$ws.Range("S2").Value = 0.01187374399721234
$ws.Range("T2").Value = 0.01366897704702647
$ws.Range("G3").Value = 3.842093
$ws.Range("H3").Value = 11.526279
$ws.Range("I3").Value = 0.8773182083458525
$ws.Range("J3").Value = 0.886145240040679
$ws.Range("O3").Value = 0.6185519418990597
$ws.Range("P3").Value = 0.704979911415303
$ws.Range("Q3").Value = 177.3270407810753
$ws.Range("R3").Value = 1595.943367029678
$ws.Range("S3").Value = 0.5426668814357309
$ws.Range("T3").Value = 0.6247145928249703
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 3.842093
$ws.Range("H4").Value = 11.526279
$ws.Range("I4").Value = 0.8773182083458525
$ws.Range("J4").Value = 0.886145240040679
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.009315666666666667
$ws.Range("N4").Value = 0.027947
$ws.Range("O4").Value = 0.0001248484115599408
$ws.Range("P4").Value = 0.000142293017222847
$ws.Range("Q4").Value = 0.03579165769033334
$ws.Range("R4").Value = 0.322124919213
$ws.Range("S4").Value = 0.0001095317847445929
$ws.Range("T4").Value = 0.0001260922799030522
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 3.842093
$ws.Range("H5").Value = 11.526279
$ws.Range("I5").Value = 0.8773182083458525
$ws.Range("J5").Value = 0.886145240040679
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 27.4428835
$ws.Range("N5").Value = 54.885767
$ws.Range("O5").Value = 0.3677890736321797
$ws.Range("P5").Value = 0.2794525848577725
$ws.Range("Q5").Value = 105.4381105951655
$ws.Range("R5").Value = 632.628663570993
$ws.Range("S5").Value = 0.3226680511281647
$ws.Range("T5").Value = 0.247635577888779
$ws.Range("I6").Value = 0.06266066604468346
$ws.Range("J6").Value = 0.06329111880393738
$ws.Range("M6").Value = 1.009860666666667
$ws.Range("N6").Value = 3.029582
$ws.Range("O6").Value = 0.01353413605720072
$ws.Range("P6").Value = 0.01542521070970148
$ws.Range("Q6").Value = 0.2771195683624444
$ws.Range("R6").Value = 2.494076115262
$ws.Range("S6").Value = 0.0008480579796835633
$ws.Range("T6").Value = 0.0009762788436034836
$ws.Range("I7").Value = 0.06266066604468346
$ws.Range("J7").Value = 0.06329111880393738
$ws.Range("O7").Value = 0.6185519418990597
$ws.Range("P7").Value = 0.704979911415303
$ws.Range("S7").Value = 0.03875887666262742
$ws.Range("T7").Value = 0.0446189673277752
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("I8").Value = 0.06266066604468346
$ws.Range("J8").Value = 0.06329111880393738
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.009315666666666667
$ws.Range("N8").Value = 0.027947
$ws.Range("O8").Value = 0.0001248484115599408
$ws.Range("P8").Value = 0.000142293017222847
$ws.Range("Q8").Value = 0.002556346247444444
$ws.Range("R8").Value = 0.023007116227
$ws.Range("S8").Value = 0.000007823084622966648
$ws.Range("T8").Value = 0.000009005884258021917
$ws.Range("D9").Value = "MuSCs"
$ws.Range("I9").Value = 0.06266066604468346
$ws.Range("J9").Value = 0.06329111880393738
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 27.4428835
$ws.Range("N9").Value = 54.885767
$ws.Range("O9").Value = 0.3677890736321797
$ws.Range("P9").Value = 0.2794525848577725
$ws.Range("Q9").Value = 7.530702285141166
$ws.Range("R9").Value = 45.184213710847
$ws.Range("S9").Value = 0.0230459083177495
$ws.Range("T9").Value = 0.01768686674830067
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.006443333333333333
$ws.Range("H10").Value = 0.01933
$ws.Range("I10").Value = 0.001471295373582865
$ws.Range("J10").Value = 0.001486098635126421
$ws.Range("M10").Value = 1.009860666666667
$ws.Range("N10").Value = 3.029582
$ws.Range("O10").Value = 0.01353413605720072
$ws.Range("P10").Value = 0.01542521070970148
$ws.Range("Q10").Value = 0.006506868895555556
$ws.Range("R10").Value = 0.05856182006
$ws.Range("S10").Value = 0.00001991271176640046
$ws.Range("T10").Value = 0.00002292338458222481
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.006443333333333333
$ws.Range("H11").Value = 0.01933
$ws.Range("I11").Value = 0.001471295373582865
$ws.Range("J11").Value = 0.001486098635126421
$ws.Range("O11").Value = 0.6185519418990597
$ws.Range("P11").Value = 0.704979911415303
$ws.Range("Q11").Value = 0.2973840645622222
$ws.Range("R11").Value = 2.67645658106
$ws.Range("S11").Value = 0.0009100726104367834
$ws.Range("T11").Value = 0.001047669684145827
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.006443333333333333
$ws.Range("H12").Value = 0.01933
$ws.Range("I12").Value = 0.001471295373582865
$ws.Range("J12").Value = 0.001486098635126421
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.009315666666666667
$ws.Range("N12").Value = 0.027947
$ws.Range("O12").Value = 0.0001248484115599408
$ws.Range("P12").Value = 0.000142293017222847
$ws.Range("Q12").Value = 0.00006002394555555555
$ws.Range("R12").Value = 0.0005402155099999999
$ws.Range("S12").Value = 0.0000001836888903273103
$ws.Range("T12").Value = 0.0000002114614586828932
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.006443333333333333
$ws.Range("H13").Value = 0.01933
$ws.Range("I13").Value = 0.001471295373582865
$ws.Range("J13").Value = 0.001486098635126421
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 27.4428835
$ws.Range("N13").Value = 54.885767
$ws.Range("O13").Value = 0.3677890736321797
$ws.Range("P13").Value = 0.2794525848577725
$ws.Range("Q13").Value = 0.1768236460183333
$ws.Range("R13").Value = 1.06094187611
$ws.Range("S13").Value = 0.0005411263624893535
$ws.Range("T13").Value = 0.000415294104939686
$ws.Range("G14").Value = 0.1308705
$ws.Range("H14").Value = 0.261741
$ws.Range("I14").Value = 0.02988347043794252
$ws.Range("J14").Value = 0.02012275958906489
$ws.Range("M14").Value = 1.009860666666667
$ws.Range("N14").Value = 3.029582
$ws.Range("O14").Value = 0.01353413605720072
$ws.Range("P14").Value = 0.01542521070970148
$ws.Range("Q14").Value = 0.132160970377
$ws.Range("R14").Value = 0.792965822262
$ws.Range("S14").Value = 0.0004044469547684497
$ws.Range("T14").Value = 0.000310397806721992
$ws.Range("G15").Value = 0.1308705
$ws.Range("H15").Value = 0.261741
$ws.Range("I15").Value = 0.02988347043794252
$ws.Range("J15").Value = 0.02012275958906489
$ws.Range("O15").Value = 0.6185519418990597
$ws.Range("P15").Value = 0.704979911415303
$ws.Range("Q15").Value = 6.040165735326999
$ws.Range("R15").Value = 36.24099441196199
$ws.Range("S15").Value = 0.01848447867007249
$ws.Range("T15").Value = 0.01418614127253041
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("G16").Value = 0.1308705
$ws.Range("H16").Value = 0.261741
$ws.Range("I16").Value = 0.02988347043794252
$ws.Range("J16").Value = 0.02012275958906489
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.009315666666666667
$ws.Range("N16").Value = 0.027947
$ws.Range("O16").Value = 0.0001248484115599408
$ws.Range("P16").Value = 0.000142293017222847
$ws.Range("Q16").Value = 0.0012191459545
$ws.Range("R16").Value = 0.007314875727
$ws.Range("S16").Value = 0.000003730903816075572
$ws.Range("T16").Value = 0.00000286332817677802
$ws.Range("D17").Value = "MuSCs"
$ws.Range("G17").Value = 0.1308705
$ws.Range("H17").Value = 0.261741
$ws.Range("I17").Value = 0.02988347043794252
$ws.Range("J17").Value = 0.02012275958906489
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 27.4428835
$ws.Range("N17").Value = 54.885767
$ws.Range("O17").Value = 0.3677890736321797
$ws.Range("P17").Value = 0.2794525848577725
$ws.Range("Q17").Value = 3.59146388508675
$ws.Range("R17").Value = 14.365855540347
$ws.Range("S17").Value = 0.0109908139092855
$ws.Range("T17").Value = 0.005623357181635713
$ws.Range("G18").Value = 0.1255403333333333
$ws.Range("H18").Value = 0.376621
$ws.Range("I18").Value = 0.02866635979793855
$ws.Range("J18").Value = 0.02895478293119232
$ws.Range("M18").Value = 1.009860666666667
$ws.Range("N18").Value = 3.029582
$ws.Range("O18").Value = 0.01353413605720072
$ws.Range("P18").Value = 0.01542521070970148
$ws.Range("Q18").Value = 0.1267782447135556
$ws.Range("R18").Value = 1.141004202422
$ws.Range("S18").Value = 0.0003879744137699693
$ws.Range("T18").Value = 0.0004466336277673094
$ws.Range("G19").Value = 0.1255403333333333
$ws.Range("H19").Value = 0.376621
$ws.Range("I19").Value = 0.02866635979793855
$ws.Range("J19").Value = 0.02895478293119232
$ws.Range("O19").Value = 0.6185519418990597
$ws.Range("P19").Value = 0.704979911415303
$ws.Range("Q19").Value = 5.794158498680222
$ws.Range("R19").Value = 52.14742648812199
$ws.Range("S19").Value = 0.01773163252019203
$ws.Range("T19").Value = 0.02041254030588129
$ws.Range("D20").Value = "Inflammatory-Mac"
$ws.Range("G20").Value = 0.1255403333333333
$ws.Range("H20").Value = 0.376621
$ws.Range("I20").Value = 0.02866635979793855
$ws.Range("J20").Value = 0.02895478293119232
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.3333333333333333
$ws.Range("M20").Value = 0.009315666666666667
$ws.Range("N20").Value = 0.027947
$ws.Range("O20").Value = 0.0001248484115599408
$ws.Range("P20").Value = 0.000142293017222847
$ws.Range("Q20").Value = 0.001169491898555556
$ws.Range("R20").Value = 0.010525427087
$ws.Range("S20").Value = 0.000003578949485978373
$ws.Range("T20").Value = 0.000004120063426311945
$ws.Range("D21").Value = "MuSCs"
$ws.Range("G21").Value = 0.1255403333333333
$ws.Range("H21").Value = 0.376621
$ws.Range("I21").Value = 0.02866635979793855
$ws.Range("J21").Value = 0.02895478293119232
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 27.4428835
$ws.Range("N21").Value = 54.885767
$ws.Range("O21").Value = 0.3677890736321797
$ws.Range("P21").Value = 0.2794525848577725
$ws.Range("Q21").Value = 3.445188742217833
$ws.Range("R21").Value = 20.671132453307
$ws.Range("S21").Value = 0.01054317391449057
$ws.Range("T21").Value = 0.008091488934117406
